$d = $word.ActiveDocument

$d.Content.Find.Execute("552÷6=92, 0", $true, $false, $false, $false, $false, $true, 1, $false, "661÷7=94, 3", 2) | Out-Null
$d.Content.Find.Execute("633÷9=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "354÷2=177, 0", 2) | Out-Null
$d.Content.Find.Execute("844÷4=211, 0", $true, $false, $false, $false, $false, $true, 1, $false, "171÷7=24, 3", 2) | Out-Null
$d.Content.Find.Execute("394÷3=131, 1", $true, $false, $false, $false, $false, $true, 1, $false, "424÷4=106, 0", 2) | Out-Null
$d.Content.Find.Execute("118÷8=14, 6", $true, $false, $false, $false, $false, $true, 1, $false, "683÷7=97, 4", 2) | Out-Null
$d.Content.Find.Execute("735÷7=105, 0", $true, $false, $false, $false, $false, $true, 1, $false, "171÷8=21, 3", 2) | Out-Null
$d.Content.Find.Execute("219÷9=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "650÷2=325, 0", 2) | Out-Null
$d.Content.Find.Execute("507÷7=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "184÷3=61, 1", 2) | Out-Null
$d.Content.Find.Execute("581÷3=193, 2", $true, $false, $false, $false, $false, $true, 1, $false, "236÷6=39, 2", 2) | Out-Null
$d.Content.Find.Execute("698÷2=349, 0", $true, $false, $false, $false, $false, $true, 1, $false, "149÷3=49, 2", 2) | Out-Null
$d.Content.Find.Execute("413÷3=137, 2", $true, $false, $false, $false, $false, $true, 1, $false, "987÷4=246, 3", 2) | Out-Null
$d.Content.Find.Execute("282÷7=40, 2", $true, $false, $false, $false, $false, $true, 1, $false, "800÷7=114, 2", 2) | Out-Null
$d.Content.Find.Execute("621÷4=155, 1", $true, $false, $false, $false, $false, $true, 1, $false, "997÷9=110, 7", 2) | Out-Null
$d.Content.Find.Execute("664÷2=332, 0", $true, $false, $false, $false, $false, $true, 1, $false, "858÷2=429, 0", 2) | Out-Null
$d.Content.Find.Execute("611÷6=101, 5", $true, $false, $false, $false, $false, $true, 1, $false, "584÷7=83, 3", 2) | Out-Null
$d.Content.Find.Execute("350÷3=116, 2", $true, $false, $false, $false, $false, $true, 1, $false, "869÷8=108, 5", 2) | Out-Null
$d.Content.Find.Execute("458÷4=114, 2", $true, $false, $false, $false, $false, $true, 1, $false, "223÷6=37, 1", 2) | Out-Null
$d.Content.Find.Execute("551÷3=183, 2", $true, $false, $false, $false, $false, $true, 1, $false, "348÷8=43, 4", 2) | Out-Null
$d.Content.Find.Execute("732÷6=122, 0", $true, $false, $false, $false, $false, $true, 1, $false, "442÷4=110, 2", 2) | Out-Null
$d.Content.Find.Execute("765÷2=382, 1", $true, $false, $false, $false, $false, $true, 1, $false, "404÷9=44, 8", 2) | Out-Null
$d.Content.Find.Execute("839÷2=419, 1", $true, $false, $false, $false, $false, $true, 1, $false, "764÷2=382, 0", 2) | Out-Null
$d.Content.Find.Execute("808÷7=115, 3", $true, $false, $false, $false, $false, $true, 1, $false, "659÷7=94, 1", 2) | Out-Null
$d.Content.Find.Execute("812÷6=135, 2", $true, $false, $false, $false, $false, $true, 1, $false, "987÷8=123, 3", 2) | Out-Null
$d.Content.Find.Execute("903÷9=100, 3", $true, $false, $false, $false, $false, $true, 1, $false, "736÷9=81, 7", 2) | Out-Null
$d.Content.Find.Execute("441÷2=220, 1", $true, $false, $false, $false, $false, $true, 1, $false, "202÷3=67, 1", 2) | Out-Null
